$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting from the existing header cell (H1) onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I (I0) and J (IF), rows 2-14
$data = @(
    @(6, 8),
    @(7, 7),
    @(7, 9),
    @(8, 8),
    @(5, 5),
    @(5, 7),
    @(3, 4),
    @(7, 7),
    @(6, 8),
    @(1, 3),
    @(1, 3),
    @(4, 5),
    @(3, 3)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
